$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.139414072036743
$ws.Range("B1").Value = 3.937014102935791
$ws.Range("C1").Value = 1.617393374443054
$ws.Range("D1").Value = 0.816184937953949
$ws.Range("E1").Value = 0.8452682495117188
